# daily auto push: 2026-01-31 18:49 UTC
#
# The source data table (A:D = 日付/曜日/時刻/ランキング) gains one new
# record. A new row is inserted at row 732 (pushing the former rows
# 732-773 down to 733-774) and is populated with:
#   A732 = 2026/01/31   (text, not a date serial)
#   B732 = 土
#   C732 = 22
#   D732 = 22
# The sheet's used range / dimension grows from D773 to D774.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 732, shifting rows 732:773 down to 733:774.
$ws.Rows.Item(732).Insert()

# Column A holds dates stored as plain text (e.g. "2026/12/29"), so force
# text formatting before assignment to avoid Excel auto-converting the
# string into a date serial number, then restore the default style so the
# new cell matches the formatting of the rest of the table.
$ws.Range("A732").NumberFormat = "@"
$ws.Range("A732").Value = "2026/01/31"
$ws.Range("A732").Style = "Normal"

$ws.Range("B732").Value = "土"
$ws.Range("C732").Value = 22
$ws.Range("D732").Value = 22
